# Add ability to load external files into TPA
#
# 1. Record a new ZEXALL compliance test run as row 20 of the "Compliance"
#    table (Table1), extending the table from A6:F19 to A6:F20.
# 2. Update the scratch "Worksheet" sheet (the TPA entry panel) with the
#    new test's raw memory dump / flag values, replacing the stale DAA/CCF
#    debugging values that were there before.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Compliance sheet: append a new test-run row to Table1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Compliance")
$tbl = $ws1.ListObjects.Item("Table1")

# Grow the table by one row (this also updates the table ref + autoFilter).
$newRow = $tbl.ListRows.Add()

# Copy the formatting straight down from the previous last row so the new
# row picks up the same number formats (date, comma, percent) without
# touching column F ("Notable fixes"), which stays blank for this entry.
$ws1.Range("A19:E19").Copy()
$ws1.Range("A20:E20").PasteSpecial(-4122)

$ws1.Range("A20").Value = 45118.736805555556
$ws1.Range("B20").Value = 1446235
$ws1.Range("C20").Value = 8326
$ws1.Range("D20").Formula = "=Table1[[#This Row],[Failures     ]]/Table1[[#This Row],[Tests     ]]"
$ws1.Range("E20").Formula = "=(`$C`$7-Table1[[#This Row],[Failures     ]])/`$C`$7"

$ws1.Range("A21").Select() | Out-Null

# ---------------------------------------------------------------------
# Worksheet sheet: refresh the TPA scratch/debug values
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Worksheet")

$ws2.Range("C2").Value = "4121FA09601D59A55B8D7990009A9D29"
$ws2.Range("C4").Value = "4121FA09601D59A55B8D799094A09D29"
$ws2.Range("C3").Value = "4121FA09601D59A55B8D799055009D29"
$ws2.Range("D3").Value = "'55"
$ws2.Range("D4").Value = "'94"
$ws2.Range("C6").Value = "A: 9A -> 00"
$ws2.Range("E3").Value = "_ Z _ H _ P _ C"
$ws2.Range("E4").Value = "S _ _ H _ P _ _"

$ws2.Range("A7").Select() | Out-Null

$ws1.Activate() | Out-Null
